$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.949.48'
$ws.Range("E2").Value = '  +1.17%  '

# Row 3
$ws.Range("D3").Value = '1.641.09'
$ws.Range("E3").Value = '  +0.41%  '

# Row 4
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
$ws.Range("E5").Value = '  +0.52%  '

# Row 6
$ws.Range("E6").Value = '  +0.66%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.72'
$ws.Range("E8").Value = '  +1.70%  '

# Row 9
$ws.Range("E9").Value = '  -0.72%  '

# Row 10
$ws.Range("E10").Value = '  +0.52%  '

# Row 11
$ws.Range("E11").Value = '  +0.88%  '

# Row 12
$ws.Range("D12").Value = '1.873.33'
$ws.Range("E12").Value = '  +0.40%  '

# Row 13
$ws.Range("D13").Value = '1.635.46'
$ws.Range("E13").Value = '  -0.02%  '

# Row 14
$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.578'
$ws.Range("E14").Value = '  +4.58%  '

# Row 15
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.10'
$ws.Range("E15").Value = '  +1.25%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.99'
$ws.Range("E16").Value = '  +1.14%  '

# Row 17
$ws.Range("D17").Value = '27.934.23'
$ws.Range("E17").Value = '  +1.17%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '231.68'
$ws.Range("E18").Value = '  +0.35%  '

# Row 19
$ws.Range("E19").Value = '  +0.58%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.62'
$ws.Range("E20").Value = '  +0.55%  '

# Row 21
$ws.Range("E21").Value = '  +0.01%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.80'
$ws.Range("E22").Value = '  +1.62%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.37'
$ws.Range("E23").Value = '  +0.30%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.08'
$ws.Range("E24").Value = '  -2.82%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.53'
$ws.Range("E25").Value = '  +1.27%  '

# Row 26
$ws.Range("E26").Value = '  +0.88%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.71'
$ws.Range("E27").Value = '  +0.96%  '

# Row 28
$ws.Range("E28").Value = '  -0.02%  '

# Row 29
$ws.Range("E29").Value = '  -0.02%  '

# Row 30
$ws.Range("E30").Value = '  +0.67%  '

# Row 31
$ws.Range("E31").Value = '  +0.09%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.34'
$ws.Range("E32").Value = '  +1.82%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.11'
$ws.Range("E33").Value = '  +0.57%  '

# Row 34
$ws.Range("D34").Value = '1.401.74'
$ws.Range("E34").Value = '  -5.19%  '

# Row 35
$ws.Range("E35").Value = '  +1.75%  '

# Row 36
$ws.Range("E36").Value = '  +0.59%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.890'
$ws.Range("E37").Value = '  +1.18%  '

# Row 38
$ws.Range("E38").Value = '  +0.45%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.922'
$ws.Range("E39").Value = '  -2.04%  '

# Row 40
$ws.Range("E40").Value = '  -0.80%  '

# Row 41
$ws.Range("E41").Value = '  -0.68%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  -0.04%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.85'
$ws.Range("E43").Value = '  +6.04%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '66.29'
$ws.Range("E44").Value = '  -2.55%  '

# Row 45
$ws.Range("E45").Value = '  +1.63%  '

# Row 46
$ws.Range("E46").Value = '  +0.04%  '

# Row 47
$ws.Range("D47").Value = '1.781.81'
$ws.Range("E47").Value = '  +0.46%  '

# Row 48
$ws.Range("E48").Value = '  +0.42%  '

# Row 49
$ws.Range("E49").Value = '  +1.05%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0506'
$ws.Range("E50").Value = '  +0.32%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.63'
$ws.Range("E51").Value = '  -0.74%  '
